$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data cells keep their original text (General/string) representation
# rather than being auto-converted to numbers by Excel when values look numeric.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "57.509.02"
$ws.Range("E2").Value = "  -4.26%  "
$ws.Range("D3").Value = "2.949.33"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "555.38"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("D6").Value = "132.34"
$ws.Range("E6").Value = "  +5.78%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +3.11%  "
$ws.Range("D9").Value = "2.942.33"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").Value = "4.81"
$ws.Range("E11").Value = "  -5.18%  "
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "32.80"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "3.433.79"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("E17").Value = "  +10.22%  "
$ws.Range("D18").Value = "2.940.42"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "57.510.09"
$ws.Range("E19").Value = "  -4.12%  "
$ws.Range("D20").Value = "417.82"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "0.683"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").Value = "6.96"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").Value = "79.26"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "7.55"
$ws.Range("E29").Value = "  +4.70%  "
$ws.Range("D30").Value = "1.98"
$ws.Range("E30").Value = "  +5.68%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "6.08"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "25.16"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +9.35%  "
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "0.936"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "2.10"
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "0.0₃0685"
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("D39").Value = "8.50"
$ws.Range("E39").Value = "  +7.18%  "
$ws.Range("D40").Value = "2.57"
$ws.Range("E40").Value = "  +4.25%  "
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "379.29"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0348"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("D44").Value = "2.654.79"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "0.241"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("D47").Value = "122.75"
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "23.42"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("E51").Value = "  +0.53%  "
